# Saldo.xlsx update — refresh the "Export" data dump with a new extract:
#   - account 005142624 (RODRIGO) no longer appears and is removed
#   - five accounts have updated Saldo (balance) figures
#   - the sheet stays sorted by Saldo, descending (as it was before),
#     and the trailing blank row + "Filtros aplicados" note are left alone

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Walk down column A from row 2 to find the last *contiguous* data row —
# the table is followed by a blank row and then a filter-notes row, so we
# stop as soon as we hit the first blank cell rather than trusting
# End(xlUp), which would land on the non-blank filter-notes row instead.
$lastDataRow = 1
$r = 2
while (-not [string]::IsNullOrEmpty($ws.Cells.Item($r, 1).Value2)) {
    $lastDataRow = $r
    $r = $r + 1
}

# New Saldo values, keyed by account number (column A).
$updates = @{
    "005646524" = 109.74    # EVANGELINA
    "000806386" = 492.08    # FERNANDA
    "004500087" = 3387.29   # JANNE
    "005331090" = 2340.86   # BEATRIZ
    "005295509" = 99.54     # BHRUNA
}
$removeAccount = "005142624"   # RODRIGO, dropped from this extract

$deleteRow = 0

for ($row = 2; $row -le $lastDataRow; $row++) {
    $conta = "{0}" -f $ws.Cells.Item($row, 1).Value2

    if ($updates.ContainsKey($conta)) {
        $ws.Cells.Item($row, 3).Value = $updates[$conta]
    }
    elseif ($conta -eq $removeAccount) {
        $deleteRow = $row
    }
}

if ($deleteRow -gt 0) {
    $ws.Rows.Item($deleteRow).Delete()
    $lastDataRow = $lastDataRow - 1
}

# Re-sort just the data block (A2:C<lastDataRow>) by Saldo (column C),
# descending — matching the sheet's existing sort order — without
# touching the blank/notes rows that follow.
$sortRange = $ws.Range("A2:C" + $lastDataRow)
$sortRange.Sort($ws.Range("C2"), 2, $null, $null, $null, $null, $null, 0)
